# Update cryptocurrency Price (column D) and Volume(1h) (column E) values
# on the symbol list sheet, matching the GitHub Actions data refresh.
# Values are stored as text (not numbers/percentages), so each is written
# with a leading apostrophe to force text entry, then the cell style is
# reset to match its row's untouched "F" cell (default style) so no
# unintended numeric/percentage formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'278.01"
$ws.Range("D2").Style = $ws.Range("F2").Style
$ws.Range("E2").Value = "'6.55%"
$ws.Range("E2").Style = $ws.Range("F2").Style
$ws.Range("D3").Value = "'27.39"
$ws.Range("D3").Style = $ws.Range("F3").Style
$ws.Range("E3").Value = "'2.00%"
$ws.Range("E3").Style = $ws.Range("F3").Style
$ws.Range("D4").Value = "'4.779"
$ws.Range("D4").Style = $ws.Range("F4").Style
$ws.Range("E4").Value = "'1.48%"
$ws.Range("E4").Style = $ws.Range("F4").Style
$ws.Range("D5").Value = "'0.06256"
$ws.Range("D5").Style = $ws.Range("F5").Style
$ws.Range("E5").Value = "'0.57%"
$ws.Range("E5").Style = $ws.Range("F5").Style
$ws.Range("D6").Value = "'6.922"
$ws.Range("D6").Style = $ws.Range("F6").Style
$ws.Range("E6").Value = "'2.66%"
$ws.Range("E6").Style = $ws.Range("F6").Style
$ws.Range("D7").Value = "'0.8807"
$ws.Range("D7").Style = $ws.Range("F7").Style
$ws.Range("E7").Value = "'3.51%"
$ws.Range("E7").Style = $ws.Range("F7").Style
$ws.Range("E8").Value = "'3.01%"
$ws.Range("E8").Style = $ws.Range("F8").Style
$ws.Range("D9").Value = "'0.1452"
$ws.Range("D9").Style = $ws.Range("F9").Style
$ws.Range("E9").Value = "'3.51%"
$ws.Range("E9").Style = $ws.Range("F9").Style
$ws.Range("D10").Value = "'0.05346"
$ws.Range("D10").Style = $ws.Range("F10").Style
$ws.Range("E10").Value = "'7.98%"
$ws.Range("E10").Style = $ws.Range("F10").Style
$ws.Range("D11").Value = "'0.07345"
$ws.Range("D11").Style = $ws.Range("F11").Style
$ws.Range("E11").Value = "'3.72%"
$ws.Range("E11").Style = $ws.Range("F11").Style
$ws.Range("D12").Value = "'0.03119"
$ws.Range("D12").Style = $ws.Range("F12").Style
$ws.Range("E12").Value = "'0.73%"
$ws.Range("E12").Style = $ws.Range("F12").Style
$ws.Range("D13").Value = "'0.09061"
$ws.Range("D13").Style = $ws.Range("F13").Style
$ws.Range("E13").Value = "'0.12%"
$ws.Range("E13").Style = $ws.Range("F13").Style
$ws.Range("E14").Value = "'1.51%"
$ws.Range("E14").Style = $ws.Range("F14").Style
$ws.Range("D15").Value = "'0.0006251"
$ws.Range("D15").Style = $ws.Range("F15").Style
$ws.Range("E15").Value = "'1.46%"
$ws.Range("E15").Style = $ws.Range("F15").Style
$ws.Range("D16").Value = "'0.005858"
$ws.Range("D16").Style = $ws.Range("F16").Style
$ws.Range("E16").Value = "'-1.86%"
$ws.Range("E16").Style = $ws.Range("F16").Style
$ws.Range("D17").Value = "'3.450"
$ws.Range("D17").Style = $ws.Range("F17").Style
$ws.Range("D18").Value = "'3.262"
$ws.Range("D18").Style = $ws.Range("F18").Style
$ws.Range("E18").Value = "'2.82%"
$ws.Range("E18").Style = $ws.Range("F18").Style
$ws.Range("E19").Value = "'5.52%"
$ws.Range("E19").Style = $ws.Range("F19").Style
$ws.Range("E21").Value = "'0.11%"
$ws.Range("E21").Style = $ws.Range("F21").Style
$ws.Range("D22").Value = "'3.856"
$ws.Range("D22").Style = $ws.Range("F22").Style
$ws.Range("E22").Value = "'-5.85%"
$ws.Range("E22").Style = $ws.Range("F22").Style
$ws.Range("E23").Value = "'1.68%"
$ws.Range("E23").Style = $ws.Range("F23").Style
$ws.Range("D24").Value = "'0.001178"
$ws.Range("D24").Style = $ws.Range("F24").Style
$ws.Range("E24").Value = "'-2.03%"
$ws.Range("E24").Style = $ws.Range("F24").Style
$ws.Range("D25").Value = "'0.004285"
$ws.Range("D25").Style = $ws.Range("F25").Style
$ws.Range("E25").Value = "'5.05%"
$ws.Range("E25").Style = $ws.Range("F25").Style
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("D26").Style = $ws.Range("F26").Style
$ws.Range("E26").Value = "'-0.02%"
$ws.Range("E26").Style = $ws.Range("F26").Style
$ws.Range("D27").Value = "'0.0001691"
$ws.Range("D27").Style = $ws.Range("F27").Style
$ws.Range("E27").Value = "'3.10%"
$ws.Range("E27").Style = $ws.Range("F27").Style
$ws.Range("D40").Value = "'0.04046"
$ws.Range("D40").Style = $ws.Range("F40").Style
$ws.Range("E40").Value = "'2.01%"
$ws.Range("E40").Style = $ws.Range("F40").Style
$ws.Range("E41").Value = "'65.55%"
$ws.Range("E41").Style = $ws.Range("F41").Style
$ws.Range("D42").Value = "'0.1154"
$ws.Range("D42").Style = $ws.Range("F42").Style
$ws.Range("E42").Value = "'3.85%"
$ws.Range("E42").Style = $ws.Range("F42").Style
$ws.Range("D43").Value = "'0.002191"
$ws.Range("D43").Style = $ws.Range("F43").Style
$ws.Range("E43").Value = "'4.26%"
$ws.Range("E43").Style = $ws.Range("F43").Style
$ws.Range("D44").Value = "'0.01290"
$ws.Range("D44").Style = $ws.Range("F44").Style
$ws.Range("E44").Value = "'-2.90%"
$ws.Range("E44").Style = $ws.Range("F44").Style
$ws.Range("D45").Value = "'0.00005084"
$ws.Range("D45").Style = $ws.Range("F45").Style
$ws.Range("E45").Value = "'-1.54%"
$ws.Range("E45").Style = $ws.Range("F45").Style
$ws.Range("E46").Value = "'-0.02%"
$ws.Range("E46").Style = $ws.Range("F46").Style
$ws.Range("E47").Value = "'840.85%"
$ws.Range("E47").Style = $ws.Range("F47").Style
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = $ws.Range("F49").Style
$ws.Range("E49").Value = "'-0.02%"
$ws.Range("E49").Style = $ws.Range("F49").Style
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("D50").Style = $ws.Range("F50").Style
$ws.Range("E50").Value = "'-0.02%"
$ws.Range("E50").Style = $ws.Range("F50").Style